$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.952.92'
$ws.Range('E2').Value = '  -3.01%  '
$ws.Range('D3').Value = '3.324.45'
$ws.Range('E3').Value = '  -5.35%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''553.92'
$ws.Range('E5').Value = '  -4.27%  '
$ws.Range('D6').Value = '''172.78'
$ws.Range('E6').Value = '  -3.09%  '
$ws.Range('D7').Value = '''0.612'
$ws.Range('E7').Value = '  -3.77%  '
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '3.315.10'
$ws.Range('E9').Value = '  -5.46%  '
$ws.Range('D10').Value = '''0.622'
$ws.Range('E10').Value = '  -1.80%  '
$ws.Range('D11').Value = '''0.161'
$ws.Range('E11').Value = '  +2.86%  '
$ws.Range('D12').Value = '''53.25'
$ws.Range('E12').Value = '  -3.78%  '
$ws.Range('D13').Value = '''0.0000271'
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('D14').Value = '''9.02'
$ws.Range('E14').Value = '  -2.21%  '
$ws.Range('D15').Value = '3.843.03'
$ws.Range('E15').Value = '  -5.81%  '
$ws.Range('D16').Value = '''18.24'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('E17').Value = '  -3.44%  '
$ws.Range('D18').Value = '3.317.86'
$ws.Range('E18').Value = '  -5.63%  '
$ws.Range('D19').Value = '''11.82'
$ws.Range('E19').Value = '  -2.17%  '
$ws.Range('D20').Value = '63.831.39'
$ws.Range('E20').Value = '  -3.19%  '
$ws.Range('D21').Value = '''0.971'
$ws.Range('E21').Value = '  -3.46%  '
$ws.Range('D22').Value = '''428.96'
$ws.Range('E22').Value = '  +3.67%  '
$ws.Range('D23').Value = '''4.62'
$ws.Range('E23').Value = '  +6.02%  '
$ws.Range('D24').Value = '''4.09'
$ws.Range('E24').Value = '  -4.11%  '
$ws.Range('D25').Value = '''84.13'
$ws.Range('E25').Value = '  -2.02%  '
$ws.Range('D26').Value = '''13.35'
$ws.Range('E26').Value = '  +2.73%  '
$ws.Range('D27').Value = '''10.62'
$ws.Range('E27').Value = '  -3.62%  '
$ws.Range('D28').Value = '''2.82'
$ws.Range('E28').Value = '  -1.25%  '
$ws.Range('D29').Value = '''8.62'
$ws.Range('E29').Value = '  -5.08%  '
$ws.Range('D30').Value = '''29.63'
$ws.Range('E30').Value = '  -2.32%  '
$ws.Range('D31').Value = '''6.61'
$ws.Range('E31').Value = '  +2.24%  '
$ws.Range('D32').Value = '''591.67'
$ws.Range('E32').Value = '  -5.26%  '
$ws.Range('D33').Value = '''11.41'
$ws.Range('E33').Value = '  -2.11%  '
$ws.Range('D34').Value = '''0.107'
$ws.Range('E34').Value = '  -3.41%  '
$ws.Range('D35').Value = '''58.15'
$ws.Range('E35').Value = '  -2.62%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  -8.92%  '
$ws.Range('D38').Value = '''35.29'
$ws.Range('E38').Value = '  -5.31%  '
$ws.Range('D39').Value = '0.0₃0746'
$ws.Range('E39').Value = '  -6.77%  '
$ws.Range('D40').Value = '''3.40'
$ws.Range('E40').Value = '  -2.33%  '
$ws.Range('D41').Value = '''0.364'
$ws.Range('E41').Value = '  -4.24%  '
$ws.Range('D42').Value = '3.097.64'
$ws.Range('E42').Value = '  -6.43%  '
$ws.Range('D43').Value = '''0.998'
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('E44').Value = '  -4.19%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '''0.0406'
$ws.Range('E45').Value = '  -2.77%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '''3.17'
$ws.Range('E46').Value = '  -2.11%  '
$ws.Range('D47').Value = '''2.43'
$ws.Range('E47').Value = '  -3.01%  '
$ws.Range('E48').Value = '  -2.33%  '
$ws.Range('D49').Value = '''2.58'
$ws.Range('E49').Value = '  -5.29%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '''8.14'
$ws.Range('E50').Value = '  -5.02%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '''132.50'
$ws.Range('E51').Value = '  -4.76%  '
